$d = $word.ActiveDocument

# Locate the "Introducción" (Heading 2) paragraph. Right after it there is a
# centered paragraph containing the screenshot image (a duplicate 404 capture,
# same hash as the ones removed from the other 9 documents) and then a
# paragraph with the descriptive text about the project. Both of those two
# paragraphs must be removed entirely, leaving the heading followed directly
# by "Desarrollo de las partes".

$introIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text -replace "[\r\a]+$", ""
    if ($txt -eq "Introducción") {
        $introIndex = $i
        break
    }
}

if ($introIndex -gt 0) {
    # Delete the paragraph right after "Introducción" (the image paragraph)
    # twice -- after the first delete, the following text paragraph shifts
    # into the same index.
    $d.Paragraphs.Item($introIndex + 1).Range.Delete()
    $d.Paragraphs.Item($introIndex + 1).Range.Delete()
}
